$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text formatting on the Price/Volume data range so numeric-looking
# strings (e.g. "1.00", "0.0400") are preserved exactly as text, matching
# the source data which stores these as inline/shared strings, not numbers.
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range("D2").Value = "63.477.02"
$ws.Range("E2").Value = "  -4.71%  "
$ws.Range("D3").Value = "3.092.48"
$ws.Range("E3").Value = "  -4.79%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "549.02"
$ws.Range("E5").Value = "  -5.51%  "
$ws.Range("D6").Value = "137.33"
$ws.Range("E6").Value = "  -10.68%  "
$ws.Range("D8").Value = "3.091.97"
$ws.Range("E8").Value = "  -4.61%  "
$ws.Range("D9").Value = "0.497"
$ws.Range("E9").Value = "  -3.30%  "
$ws.Range("E10").Value = "  -5.20%  "
$ws.Range("D11").Value = "6.33"
$ws.Range("E11").Value = "  -11.25%  "
$ws.Range("D12").Value = "0.471"
$ws.Range("E12").Value = "  -3.43%  "
$ws.Range("D13").Value = "35.48"
$ws.Range("E13").Value = "  -6.41%  "
$ws.Range("D14").Value = "0.0000217"
$ws.Range("E14").Value = "  -7.49%  "
$ws.Range("D15").Value = "3.598.64"
$ws.Range("E15").Value = "  -4.60%  "
$ws.Range("D16").Value = "63.467.83"
$ws.Range("E16").Value = "  -4.78%  "
$ws.Range("E17").Value = "  -2.93%  "
$ws.Range("D18").Value = "3.090.90"
$ws.Range("E18").Value = "  -4.93%  "
$ws.Range("E19").Value = "  -4.85%  "
$ws.Range("D20").Value = "489.09"
$ws.Range("E20").Value = "  -12.21%  "
$ws.Range("D21").Value = "13.62"
$ws.Range("E21").Value = "  -5.41%  "
$ws.Range("D22").Value = "0.719"
$ws.Range("E22").Value = "  -3.04%  "
$ws.Range("D23").Value = "7.29"
$ws.Range("E23").Value = "  -6.01%  "
$ws.Range("D24").Value = "79.13"
$ws.Range("E24").Value = "  -3.21%  "
$ws.Range("D25").Value = "12.38"
$ws.Range("E25").Value = "  -9.09%  "
$ws.Range("D26").Value = "1.00"
$ws.Range("E26").Value = "  +0.07%  "
$ws.Range("D27").Value = "8.52"
$ws.Range("E27").Value = "  -7.69%  "
$ws.Range("E28").Value = "  -6.70%  "
$ws.Range("D29").Value = "1.97"
$ws.Range("E29").Value = "  -11.57%  "
$ws.Range("E30").Value = "  -0.22%  "
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("E32").Value = "  -4.84%  "
$ws.Range("D33").Value = "2.51"
$ws.Range("E33").Value = "  -9.04%  "
$ws.Range("D34").Value = "59.03"
$ws.Range("E34").Value = "  +6.47%  "
$ws.Range("D35").Value = "507.34"
$ws.Range("E35").Value = "  -9.71%  "
$ws.Range("D36").Value = "6.06"
$ws.Range("E36").Value = "  -4.79%  "
$ws.Range("D37").Value = "5.12"
$ws.Range("E37").Value = "  -10.22%  "
$ws.Range("D38").Value = "0.0400"
$ws.Range("E38").Value = "  -12.58%  "
$ws.Range("D39").Value = "3.147.88"
$ws.Range("E39").Value = "  -0.56%  "
$ws.Range("D40").Value = "0.0803"
$ws.Range("E40").Value = "  -6.72%  "
$ws.Range("E41").Value = "  -10.12%  "
$ws.Range("D42").Value = "8.18"
$ws.Range("E42").Value = "  -4.99%  "
$ws.Range("D43").Value = "2.63"
$ws.Range("E43").Value = "  -13.26%  "
$ws.Range("D44").Value = "0.257"
$ws.Range("E44").Value = "  -6.29%  "
$ws.Range("E45").Value = "  +0.13%  "
$ws.Range("D46").Value = "25.40"
$ws.Range("E46").Value = "  -4.21%  "
$ws.Range("E47").Value = "  -10.47%  "
$ws.Range("D48").Value = "120.53"
$ws.Range("E48").Value = "  -4.33%  "
$ws.Range("D49").Value = "0.109"
$ws.Range("E49").Value = "  -3.32%  "
$ws.Range("D50").Value = "0.0₃0507"
$ws.Range("E50").Value = "  -8.78%  "

# Row 51: coin replaced (CoreDAO -> ThetaToken)
$ws.Range("B51").Value = "ThetaToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D51").Value = "2.04"
$ws.Range("E51").Value = "  -8.73%  "
